# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3256
$ws1.Range("F6").Value = 2106
$ws1.Range("F8").Value = 152
$ws1.Range("F10").Value = 1200
$ws1.Range("F12").Value = 1173
$ws1.Range("F13").Value = 96

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3256
$ws4.Range("F6").Value = 2106
$ws4.Range("F9").Value = 152
$ws4.Range("F11").Value = 1200
$ws4.Range("F13").Value = 1173
$ws4.Range("F14").Value = 96
